$d = $word.ActiveDocument

$pairs = @(
    @("11×99=", "82×28="),
    @("22×11=", "58×11="),
    @("30×25=", "93×11="),
    @("48×41=", "29×51="),
    @("22×33=", "23×64="),
    @("70×59=", "90×54="),
    @("25×61=", "21×27="),
    @("79×78=", "81×88="),
    @("95×99=", "47×29="),
    @("45×15=", "22×92="),
    @("67×27=", "24×46="),
    @("23×43=", "80×26="),
    @("94×98=", "64×70="),
    @("34×37=", "74×18="),
    @("93×81=", "62×35="),
    @("78×38=", "18×37="),
    @("89×35=", "24×35="),
    @("25×24=", "89×63="),
    @("76×15=", "55×56="),
    @("97×55=", "41×66="),
    @("43×24=", "36×40="),
    @("85×91=", "61×75="),
    @("16×96=", "41×71="),
    @("63×46=", "22×53="),
    @("43×77=", "96×99=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
